$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 158, pushing the existing row 158 (and all
# rows below it) down by one. This mirrors the diff, where every row from
# the old 158..261 reappears unchanged one row lower (159..262) and a
# brand-new record is inserted at 158.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with a new "Cilantro" record. Every
# field matches the (old) row 158 except the date (column D) and the
# volume (column J), which carry the new values from the diff.
$ws.Range("A158").Value = 8
$ws.Range("B158").Value = "Terminal La Palmera de La Serena"
$ws.Range("C158").Value = "Coquimbo"
$ws.Range("D158").Value = 45216
$ws.Range("E158").Value = 4
$ws.Range("F158").Value = 100112040
$ws.Range("G158").Value = "Cilantro"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 2000
$ws.Range("K158").Value = 1500
$ws.Range("L158").Value = 2000
$ws.Range("M158").Value = 1750
$ws.Range("N158").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O158").Value = "Provincia del Elquí"
$ws.Range("P158").Value = 1167
$ws.Range("Q158").Value = 1.5
$ws.Range("R158").Value = "Hortaliza"
